# Reorders the comma-separated training names inside the
# "missing_trainings" column (column E) so that the check for
# "Champions" trainings runs after the other checks, matching the
# updated validation order used when the report is generated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: old text -> new text (order of items within the cell changed)
$replacements = @{
    "Coffee Champions - ToT 1, HOR 1" = "HOR 1, Coffee Champions - ToT 1";
    "Coffee Champions - ToT 1, HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection, Village Champions" = "HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection, Coffee Champions - ToT 1, Village Champions";
    "CATs / Agriculture Champions, WASH Champions Training, Coffee Champions - ToT 1, Village Champions, HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection" = "WASH Champions Training, HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection, Village Champions, Coffee Champions - ToT 1, CATs / Agriculture Champions";
    "Coffee Champions - ToT 1, HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection, Village Champions, HOR 1" = "HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection, HOR 1, Coffee Champions - ToT 1, Village Champions";
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

# Find the "missing_trainings" column from the header row.
$headerRow = 1
$colCount = $usedRange.Columns.Count
$targetCol = 0
for ($c = 1; $c -le $colCount; $c++) {
    $headerValue = $ws.Cells.Item($headerRow, $c).Value2
    if ($headerValue -eq "missing_trainings") {
        $targetCol = $c
        break
    }
}

if ($targetCol -eq 0) {
    $targetCol = 5
}

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $targetCol)
    $current = $cell.Value2
    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
